$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 3
$ws.Range("H3").Value = 52828.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 52828.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 52828.5
$ws.Range("N3").Value = -53056.5

# ALC row 28
$ws.Range("H28").Value = 773.3
$ws.Range("I28").Value = 293.05264
$ws.Range("J28").Value = 9898
$ws.Range("K28").Value = 293.05264
$ws.Range("L28").Value = 9898
$ws.Range("M28").Value = 191.94736
$ws.Range("N28").Value = -10868

# ALC row 38
$ws.Range("H38").Value = 781.6667
$ws.Range("I38").Value = 382.30768
$ws.Range("J38").Value = 1820
$ws.Range("K38").Value = 1146.92304
$ws.Range("L38").Value = 5460
$ws.Range("M38").Value = -774.9230400000001
$ws.Range("N38").Value = -6204

# ALC row 41
$ws.Range("H41").Value = 241.61539
$ws.Range("I41").Value = 144.28572
$ws.Range("J41").Value = 355.16666
$ws.Range("K41").Value = 144.28572
$ws.Range("L41").Value = 355.16666
$ws.Range("M41").Value = 295.71428
$ws.Range("N41").Value = -1235.16666

# ALC row 62
$ws.Range("H62").Value = 7739.2144
$ws.Range("I62").Value = 8480
$ws.Range("J62").Value = 3294.5
$ws.Range("K62").Value = 8480
$ws.Range("L62").Value = 3294.5
$ws.Range("M62").Value = -7856
$ws.Range("N62").Value = -4542.5

# ALC row 65
$ws.Range("H65").Value = 7739.2144
$ws.Range("I65").Value = 8480
$ws.Range("J65").Value = 3294.5
$ws.Range("K65").Value = 42400
$ws.Range("L65").Value = 16472.5
$ws.Range("M65").Value = -39280
$ws.Range("N65").Value = -22712.5

# ALC row 102
$ws.Range("H102").Value = 52828.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 52828.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 52828.5
$ws.Range("N102").Value = -59318.5

# ALC row 106
$ws.Range("H106").Value = 27506.047
$ws.Range("I106").Value = 5420.5713
$ws.Range("J106").Value = 38548.785
$ws.Range("K106").Value = 5420.5713
$ws.Range("L106").Value = 38548.785
$ws.Range("M106").Value = -4789.5713
$ws.Range("N106").Value = -39810.785

# ALC row 137
$ws.Range("H137").Value = 12970.18
$ws.Range("I137").Value = 5798.4736
$ws.Range("J137").Value = 19783.3
$ws.Range("K137").Value = 17395.4208
$ws.Range("L137").Value = 59349.89999999999
$ws.Range("M137").Value = -14845.4208
$ws.Range("N137").Value = -64449.89999999999

# ALC row 138
$ws.Range("H138").Value = 7256.3335
$ws.Range("I138").Value = 6360.625
$ws.Range("J138").Value = 7422.9766
$ws.Range("K138").Value = 19081.875
$ws.Range("L138").Value = 22268.9298
$ws.Range("M138").Value = -13941.875
$ws.Range("N138").Value = -32548.9298


$ws = $wb.Worksheets.Item("ARM")
# ARM row 6
$ws.Range("H6").Value = 18666.666
$ws.Range("I6").Value = 42000
$ws.Range("J6").Value = 7000
$ws.Range("K6").Value = 42000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = -41827
$ws.Range("N6").Value = -7346

# ARM row 32
$ws.Range("H32").Value = 17433.762
$ws.Range("I32").Value = 14788.368
$ws.Range("J32").Value = 29999.375
$ws.Range("K32").Value = 14788.368
$ws.Range("L32").Value = 29999.375
$ws.Range("M32").Value = -14501.368
$ws.Range("N32").Value = -30573.375

# ARM row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0

# ARM row 61
$ws.Range("H61").Value = 11660.066
$ws.Range("I61").Value = 3667.3333
$ws.Range("J61").Value = 23649.166
$ws.Range("K61").Value = 3667.3333
$ws.Range("L61").Value = 23649.166
$ws.Range("M61").Value = -3455.3333
$ws.Range("N61").Value = -24073.166

# ARM row 110
$ws.Range("H110").Value = 2118.3914
$ws.Range("I110").Value = 1942
$ws.Range("J110").Value = 5999
$ws.Range("K110").Value = 1942
$ws.Range("L110").Value = 5999
$ws.Range("M110").Value = 103
$ws.Range("N110").Value = -10089

# ARM row 122
$ws.Range("H122").Value = 3874.862
$ws.Range("I122").Value = 3422.4375
$ws.Range("J122").Value = 4431.6924
$ws.Range("K122").Value = 10267.3125
$ws.Range("L122").Value = 13295.0772
$ws.Range("M122").Value = -7817.3125
$ws.Range("N122").Value = -18195.0772

# ARM row 136
$ws.Range("H136").Value = 11660.066
$ws.Range("I136").Value = 3667.3333
$ws.Range("J136").Value = 23649.166
$ws.Range("K136").Value = 11001.9999
$ws.Range("L136").Value = 70947.49800000001
$ws.Range("M136").Value = -8451.999899999999
$ws.Range("N136").Value = -76047.49800000001


$ws = $wb.Worksheets.Item("BSM")
# BSM row 37
$ws.Range("H37").Value = 1065
$ws.Range("I37").Value = 509
$ws.Range("J37").Value = 2733
$ws.Range("K37").Value = 509
$ws.Range("L37").Value = 2733
$ws.Range("M37").Value = -372
$ws.Range("N37").Value = -3007

# BSM row 75
$ws.Range("H75").Value = 10583.333
$ws.Range("I75").Value = 10583.333
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 10583.333
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -9647.333000000001

# BSM row 78
$ws.Range("H78").Value = 10583.333
$ws.Range("I78").Value = 10583.333
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 31749.999
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -27069.999

# BSM row 86
$ws.Range("H86").Value = 773524.6
$ws.Range("I86").Value = 3335295
$ws.Range("J86").Value = 4993.5
$ws.Range("K86").Value = 3335295
$ws.Range("L86").Value = 4993.5
$ws.Range("M86").Value = -3334172
$ws.Range("N86").Value = -7239.5

# BSM row 89
$ws.Range("H89").Value = 773524.6
$ws.Range("I89").Value = 3335295
$ws.Range("J89").Value = 4993.5
$ws.Range("K89").Value = 16676475
$ws.Range("L89").Value = 24967.5
$ws.Range("M89").Value = -16670859
$ws.Range("N89").Value = -36199.5

# BSM row 134
$ws.Range("H134").Value = 9350.609
$ws.Range("I134").Value = 4365.1377
$ws.Range("J134").Value = 21398.834
$ws.Range("K134").Value = 13095.4131
$ws.Range("L134").Value = 64196.50199999999
$ws.Range("M134").Value = -10560.4131
$ws.Range("N134").Value = -69266.50199999999


$ws = $wb.Worksheets.Item("CRP")
# CRP row 2
$ws.Range("H2").Value = 4960.3335
$ws.Range("I2").Value = 1952.6
$ws.Range("J2").Value = 19999
$ws.Range("K2").Value = 1952.6
$ws.Range("L2").Value = 19999
$ws.Range("M2").Value = -1839.6
$ws.Range("N2").Value = -20225

# CRP row 31
$ws.Range("H31").Value = 4112.8184
$ws.Range("I31").Value = 2515
$ws.Range("J31").Value = 5025.857
$ws.Range("K31").Value = 2515
$ws.Range("L31").Value = 5025.857
$ws.Range("M31").Value = -2220
$ws.Range("N31").Value = -5615.857

# CRP row 34
$ws.Range("H34").Value = 4112.8184
$ws.Range("I34").Value = 2515
$ws.Range("J34").Value = 5025.857
$ws.Range("K34").Value = 2515
$ws.Range("L34").Value = 5025.857
$ws.Range("M34").Value = -2313
$ws.Range("N34").Value = -5429.857

# CRP row 133
$ws.Range("H133").Value = 1134022.2
$ws.Range("I133").Value = 70000
$ws.Range("J133").Value = 1311359.4
$ws.Range("K133").Value = 70000
$ws.Range("L133").Value = 1311359.4
$ws.Range("M133").Value = -67470
$ws.Range("N133").Value = -1316419.4

# CRP row 141
$ws.Range("H141").Value = 259666.52
$ws.Range("I141").Value = 13358
$ws.Range("J141").Value = 388175.3
$ws.Range("K141").Value = 13358
$ws.Range("L141").Value = 388175.3
$ws.Range("M141").Value = -8178
$ws.Range("N141").Value = -398535.3


$ws = $wb.Worksheets.Item("CUL")
# CUL row 7
$ws.Range("H7").Value = 599.25
$ws.Range("I7").Value = 599.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1797.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1685.75
$ws.Range("N7").ClearContents()

# CUL row 110
$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 3000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 9000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -4910

# CUL row 114
$ws.Range("H114").Value = 608.4
$ws.Range("I114").Value = 608.4
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 1825.2
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 1428.8

# CUL row 120
$ws.Range("H120").Value = 8300
$ws.Range("I120").Value = 5750
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 17250
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -12412
$ws.Range("N120").Value = -39676

# CUL row 122
$ws.Range("H122").Value = 16667501
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 16667501
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 150007509
$ws.Range("N122").Value = -150012409

# CUL row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
# GSM row 113
$ws.Range("H113").Value = 115260.11
$ws.Range("I113").Value = 137366.8
$ws.Range("J113").Value = 4726.6665
$ws.Range("K113").Value = 137366.8
$ws.Range("L113").Value = 4726.6665
$ws.Range("M113").Value = -135196.8
$ws.Range("N113").Value = -9066.666499999999


$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Range("H46").Value = 2029.2727
$ws.Range("I46").Value = 1700
$ws.Range("J46").Value = 2044.9524
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 2044.9524
$ws.Range("M46").Value = -1512
$ws.Range("N46").Value = -2420.9524

# LTW row 99
$ws.Range("H99").Value = 47839
$ws.Range("I99").Value = 47839
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 47839
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -44844

# LTW row 132
$ws.Range("H132").Value = 7334.9536
$ws.Range("I132").Value = 6851.2354
$ws.Range("J132").Value = 9162.333000000001
$ws.Range("K132").Value = 20553.7062
$ws.Range("L132").Value = 27486.999
$ws.Range("M132").Value = -18023.7062
$ws.Range("N132").Value = -32546.999


$ws = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws.Range("H81").Value = 2798.4
$ws.Range("I81").Value = 1993
$ws.Range("J81").Value = 2999.75
$ws.Range("K81").Value = 3986
$ws.Range("L81").Value = 5999.5
$ws.Range("M81").Value = -2925
$ws.Range("N81").Value = -8121.5

# WVR row 84
$ws.Range("H84").Value = 2798.4
$ws.Range("I84").Value = 1993
$ws.Range("J84").Value = 2999.75
$ws.Range("K84").Value = 19930
$ws.Range("L84").Value = 29997.5
$ws.Range("M84").Value = -14626
$ws.Range("N84").Value = -40605.5

# WVR row 132
$ws.Range("H132").Value = 262061.39
$ws.Range("I132").Value = 458987.47
$ws.Range("J132").Value = 35596.4
$ws.Range("K132").Value = 1376962.41
$ws.Range("L132").Value = 106789.2
$ws.Range("M132").Value = -1374432.41
$ws.Range("N132").Value = -111849.2

# WVR row 136
$ws.Range("H136").Value = 7319448
$ws.Range("I136").Value = 10715693
$ws.Range("J136").Value = 4459.769
$ws.Range("K136").Value = 32147079
$ws.Range("L136").Value = 13379.307
$ws.Range("M136").Value = -32144529
$ws.Range("N136").Value = -18479.307

